$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '''26.109.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '''1.653.65'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D5").Value = '''218.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").Value = '''0.5257'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").Value = '''0.2668'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("D9").Value = '''0.06368'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").Value = '''0.07692'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '''4.598'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.52%  '
$ws.Range("D13").Value = '''1.641.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '''1.880.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("D15").Value = '''0.5609'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").Value = '''0.0₅8236'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '''65.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '''26.112.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").Value = '''4.704'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("D21").Value = '''10.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").Value = '''191.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.73%  '
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D25").Value = '''146.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("D27").Value = '''7.265'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("D28").Value = '''15.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").Value = '''1.498'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("D30").Value = '''0.05658'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").Value = '''3.504'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("D33").Value = '''3.390'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("D34").Value = '''1.582'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("D35").Value = '''2.800'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").Value = '''0.5794'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").Value = '''0.01594'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("D40").Value = '''5.981'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").Value = '''0.8414'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").Value = '''1.025.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.63%  '
$ws.Range("D44").Value = '''101.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.09%  '
$ws.Range("D45").Value = '''1.791.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").Value = '''58.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '''1.005'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").Value = '''0.05337'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.60%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.062'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.0₈103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").Value = '''0.4341'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.61%  '
